$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column AP into the new column AQ, matching the
# header-style (s=1) on row 1 and numeric-style (s=2) on rows 2-11.
$ws.Range("AP1:AP11").Copy()
$ws.Range("AQ1:AQ11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new "07-ago" column.
$ws.Range("AQ1").Value = "07-ago"
$ws.Range("AQ2").Value = 15
$ws.Range("AQ3").Value = 13
$ws.Range("AQ4").Value = 12
$ws.Range("AQ5").Value = 17
$ws.Range("AQ6").Value = 8
$ws.Range("AQ7").Value = 17
$ws.Range("AQ8").Value = 18
$ws.Range("AQ9").Value = 19
$ws.Range("AQ10").Value = 13
$ws.Range("AQ11").Value = 23

$ws.Range("AE1:AL1").EntireColumn.ColumnWidth = -0.8333333333333334
$ws.Range("AE1:AL1").EntireColumn.Hidden = $true

$ws.Range("AV8").Select() | Out-Null
